# "add % format data"
#
# - Rename the C1/D1/E1 headers and add a new F1 "百分比" (percentage) header.
# - Add a new column F holding percentage-formatted values (built-in
#   numFmtId 10, i.e. "0.00%") for the existing data rows 2-5 (row 3 / row 4
#   keep a non-numeric placeholder, matching the source data).
# - Append a new data row 6 (same A/B/C/D/E as row 5) with its own F value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
# A1 "姓名" and B1 "Total(persons)" are unchanged.
$ws.Range("C1").Value = "千分位数字"
$ws.Range("D1").Value = "Double字段"
$ws.Range("E1").Value = "数字文本混合"
$ws.Range("F1").Value = "百分比"

# --- Add new column F (percentage) for the existing data rows ---
$ws.Range("F2").Value = 0.6812
$ws.Range("F2").NumberFormat = "0.00%"

$ws.Range("F3").Value = 0.0
$ws.Range("F3").NumberFormat = "0.00%"

# Row 4's percentage cell is left as non-numeric text in the source data.
$ws.Range("F4").Value = "--"

$ws.Range("F5").Value = 0.0002
$ws.Range("F5").NumberFormat = "0.00%"

# --- Append new row 6, same A-E data as row 5, with its own F value ---
$ws.Range("A6").Value = "M军"
$ws.Range("B6").Value = "Brazil"
$ws.Range("C6").Value = 206081000
$ws.Range("C6").NumberFormat = "#,##0.00"
$ws.Range("D6").Value = 24.66
$ws.Range("E6").Value = 2016

$ws.Range("F6").Value = 1.0
$ws.Range("F6").NumberFormat = "0.00%"
